# Generate Report for Handoff
# Adds 3 new files (2 .png images + 1 .md file) that are "Ready for handoff"
# to the Overview / zh-cn / de-de sheets, growing each table from 4 to 7 rows.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f1edcdcc92880349a893c500c7f03c03f449033/e2e/"
$repoZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3d82eaa5db325258fc5e8b1fc71f28c7050fc526/e2e/"
$repoDeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8cb7b57fa3b2f77415119a4f4a5466e964c94195/e2e/"

$file1 = "090d781c-ce1e-4d89-8ee2-9336825289a8.png"
$file2 = "2fa9eaa8-46ca-497c-939b-fb5bb81a3304.png"
$file3 = "a5eeb7bb-47bd-4a7e-bcb3-5196be363ab5.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = $file1
$wsOverview.Range("B5").Value = "e2e\" + $file1
$wsOverview.Range("C5").Value = ".png"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-31 13:47:58"

$wsOverview.Range("A6").Value = $file2
$wsOverview.Range("B6").Value = "e2e\" + $file2
$wsOverview.Range("C6").Value = ".png"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-31 13:47:58"

$wsOverview.Range("A7").Value = $file3
$wsOverview.Range("B7").Value = "e2e\" + $file3
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-31 13:47:58"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $repoBase + $file1, "", "", "e2e\" + $file1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), $repoBase + $file2, "", "", "e2e\" + $file2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), $repoBase + $file3, "", "", "e2e\" + $file3)

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A5").Value = $file1
$wsZhCn.Range("B5").Value = ".png"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "60f2071b7f9a8e82f45c89e23ce67712b0225ba8.png"
$wsZhCn.Range("H5").Value = "2016-08-31 13:47:47"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M5").Value = "True(Dependency)"
$wsZhCn.Range("N5").Value = "e2e\" + $file3
$wsZhCn.Range("O5").Value = "False"

$wsZhCn.Range("A6").Value = $file2
$wsZhCn.Range("B6").Value = ".png"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "d03cb40a535177257af7062c9b5d3b28b6262a7c.png"
$wsZhCn.Range("H6").Value = "2016-08-31 13:47:47"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M6").Value = "True(Dependency)"
$wsZhCn.Range("N6").Value = "e2e\" + $file3
$wsZhCn.Range("O6").Value = "False"

$wsZhCn.Range("A7").Value = $file3
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "a5eeb7bb-47bd-4a7e-bcb3-5196be363ab5.d51ad81e8d3aa2f54e6151a95d0812fe51b25348.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-31 13:47:47"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), $repoZhCn + $file1, "", "", $file1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), $repoZhCn + $file2, "", "", $file2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), $repoZhCn + $file3, "", "", $file3)

$wsZhCn.Columns.Item(13).ColumnWidth = 16.75
$wsZhCn.Columns.Item(14).ColumnWidth = 39.09

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A5").Value = $file1
$wsDeDe.Range("B5").Value = ".png"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "60f2071b7f9a8e82f45c89e23ce67712b0225ba8.png"
$wsDeDe.Range("H5").Value = "2016-08-31 13:47:58"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M5").Value = "True(Dependency)"
$wsDeDe.Range("N5").Value = "e2e\" + $file3
$wsDeDe.Range("O5").Value = "False"

$wsDeDe.Range("A6").Value = $file2
$wsDeDe.Range("B6").Value = ".png"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "d03cb40a535177257af7062c9b5d3b28b6262a7c.png"
$wsDeDe.Range("H6").Value = "2016-08-31 13:47:58"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M6").Value = "True(Dependency)"
$wsDeDe.Range("N6").Value = "e2e\" + $file3
$wsDeDe.Range("O6").Value = "False"

$wsDeDe.Range("A7").Value = $file3
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "a5eeb7bb-47bd-4a7e-bcb3-5196be363ab5.d51ad81e8d3aa2f54e6151a95d0812fe51b25348.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-31 13:47:58"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), $repoDeDe + $file1, "", "", $file1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), $repoDeDe + $file2, "", "", $file2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), $repoDeDe + $file3, "", "", $file3)

$wsDeDe.Columns.Item(13).ColumnWidth = 16.75
$wsDeDe.Columns.Item(14).ColumnWidth = 39.09

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P7"))

Write-Host "Handoff report rows added."
